# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" across the
#    Overview sheet (zh-cn/de-de status columns E/F) and the per-locale
#    handback sheets (Status column C).
# 2. Narrow the "Status" columns to match the shorter label:
#    Overview!E:F and each locale sheet's column C.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text wherever it appears ---------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- Re-size the Status columns to fit the new (shorter) text -------------
# (target stored width ~= 13.41 "character" units; the host quantizes
# ColumnWidth writes to whole pixels, so feed it the input that lands on
# the closest achievable pixel width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
